$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")

$ws.Range("A1").Value = "Us Batch 7 1"
$ws.Range("A2").Value = "Us Batch 7 2"
$ws.Range("A3").Value = "Us Batch 7 3"
$ws.Range("A4").Value = "Us Batch 7 4"
$ws.Range("A5").Value = "Us Batch 7 5"
$ws.Range("A6").Value = "Us Batch 7 6"
$ws.Range("A7").Value = "Us Batch 7 7"
$ws.Range("A8").Value = "Us Batch 7 8"

$ws.Range("B1").Value = "ub71"
$ws.Range("B2").Value = "ub72"
$ws.Range("B3").Value = "ub73"
$ws.Range("B4").Value = "ub74"
$ws.Range("B5").Value = "ub75"
$ws.Range("B6").Value = "ub76"
$ws.Range("B7").Value = "ub77"
$ws.Range("B8").Value = "ub78"

$ws.Activate()
$ws.Range("C8").Select()

$win = $excel.ActiveWindow
$win.Zoom = 160
